$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" -----------------------------------
# Shift every Week_Start_Date (column B, rows 2-17) forward by one week.
# Row 17 previously had no data beyond what row 16 held, so it now gets a
# brand new date one week after the old row 16 value.
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

$newDates = @(
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27",
    "2025-05-04",
    "2025-05-11",
    "2025-05-18"
)

# Force the column to Text formatting first so the yyyy-mm-dd strings are
# not auto-converted into date serial numbers by Excel's smart entry.
$dateRange = $wsForecast.Range("B2:B17")
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = 2 + $i
    $wsForecast.Cells.Item($row, 2).Value = $newDates[$i]
}

# Restore default formatting so no stray number-format is left behind.
$dateRange.ClearFormats()

# --- Sheet 2: "Summary" -------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# Every value in this column is stored as plain text in the source file
# (even the numeric-looking ones), so force Text formatting before writing
# so Excel doesn't reinterpret "13" as a number or the dates as serials.
$summaryRange = $wsSummary.Range("B2:B15")
$summaryRange.NumberFormat = "@"

$wsSummary.Range("B2").Value = "2022-12-25 to 2025-01-26"
$wsSummary.Range("B4").Value = "13"
$wsSummary.Range("B8").Value = "340 units"
$wsSummary.Range("B13").Value = "2025-02-02"
$wsSummary.Range("B15").Value = "2025-02-02"

$summaryRange.ClearFormats()
